# edit.ps1
# Target edit (per commit "uploaded the daily files"): the run containing
#   "neue freunde machen oder ein neues Domain entdecken. ..."
# is split so that "freunde" is wrapped in spellStart/spellEnd proofErr
# markers (as Word's proofing pass would do once the run is re-typed) and
# "machen" is corrected to "finden":
#   "neue " + [spellStart]"freunde"[spellEnd] + " " + "finden" + " oder ein neues Domain entdecken. ..."
#
# The runtime's Range.Text / Find based editing normalizes/merges runs that
# share formatting, which loses the fine-grained run/proofErr split seen in
# the target OOXML. Range.InsertXML, when applied to the *entire* paragraph
# range, replaces just that paragraph's run content (the paragraph's own
# <w:pPr> / identity attributes are kept) with exactly the OOXML we supply -
# so we rebuild the paragraph's run sequence with the split applied and
# push it back through InsertXML.

$d = $word.ActiveDocument

# Locate the target paragraph via a stable, distinctive substring.
$anchor = $d.Content
$found = $anchor.Find.Execute("neue freunde machen oder ein neues Domain entdecken", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence in the document."
}

$targetPara = $d.Range($anchor.Start, $anchor.Start).Paragraphs(1)
$paraRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)

# Sanity-check we still have the expected original wording before mutating.
if ($paraRange.Text -notlike "*neue freunde machen oder ein neues Domain entdecken*") {
    throw "Paragraph text did not match the expected original content."
}

$newParagraphXml = '<w:p><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Das Erasmus Programm ist sehr nützlich, denn ein Student kann ein Bildungsreise machen. Das ist ein wichtiges Forschung, denn Mann kann </w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">neue </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>freunde</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>finden</w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> oder ein neues Domain entdecken. Eine Austausch Programme ist sehr nützlich. Wenn Mann </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00177350"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>ist</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00177350"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> allein im eine neues Land</w:t></w:r><w:r w:rsidR="00A64C8B"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> oder Arbeit</w:t></w:r><w:r w:rsidR="00177350"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">, hat er wenige Freuden. Mit ein Austausch Programme hat er viele Personen und er kann Freunde </w:t></w:r><w:r w:rsidR="00A64629"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>finden</w:t></w:r><w:r w:rsidR="00177350"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">. Mit ihrer neue freunden </w:t></w:r><w:r w:rsidR="005675CB"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">kann er ein Wohnung teilen. Mit ihrer freunden kann er arbeiten zusammen. Mit ihrer </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="005675CB"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>freunden</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="005675CB"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> hat er nicht oder wenige Heimat. Aber, wenn Mann </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="005675CB"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>ist</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="005675CB"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> allein, hat er viel Heimat. Dann, eine Austausch Programme ist da für ihn. Erasmus ist eine sehr gut Programme für Europa, damals Mann muss in eine Austausch Programme begleichen.</w:t></w:r><w:r w:rsidR="005675CB"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:br w:type="page"/></w:r></w:p>'

$package = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" + `
    "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" + `
    "<pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" + `
    "<w:body>" + $newParagraphXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$paraRange.InsertXML($package)

Write-Output "Paragraph updated."
